$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that changed
$updates = @{
    2  = 115
    3  = 207
    5  = 6521
    9  = 5887
    11 = 189
    14 = 82
    18 = 336
    21 = 4208
}

# Both "展览" and "全部类型" sheets received the same updates
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
